# Updated symbol list on Fri Jan 13 21:57:12 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns for the coin rows that changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.30%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.16%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.253"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.72%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07157"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.466"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.37%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.569"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.388"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.63%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9095"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.60%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1623"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.16%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07614"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "14.68%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07749"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.37%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02920"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.70%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08988"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.27%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001592"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.27%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006495"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006333"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.00%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.507"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.82%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.231"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.95%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3261"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.37%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1362"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.95%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.023"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04518"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.79%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001204"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.97%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004234"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.77%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001164"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.80%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001918"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "18.58%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04389"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.42%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006976"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.46%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1267"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.87%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002200"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.96%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01325"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.58%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005813"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.02%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01288"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-1.44%"